$d = $word.ActiveDocument

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# 1) The "Gradient descent" paragraph loses its paragraph-mark run properties
#    (the <w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr>).
#    Rewrite the whole paragraph with clean markup but identical text.
# ---------------------------------------------------------------------------
$targetText = "             Gradient descent"
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq ($targetText + "`r")) {
        $xml = '<w:p ' + $wNs + '><w:r><w:t xml:space="preserve">             Gradient descent</w:t></w:r></w:p>'
        $p.Range.InsertXML($xml)
        break
    }
}

# ---------------------------------------------------------------------------
# 2) The empty paragraph that only held the "_GoBack" bookmark (right before
#    "5.1.2 Neural Networks - Backpropagation Algorithm") becomes a plain
#    empty paragraph with no bookmark at all.
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "`r") {
        $next = $p.Next()
        if ($next -ne $null -and $next.Range.Text -eq ("5.1.2 Neural Networks - Backpropagation Algorithm" + "`r")) {
            $xml = '<w:p ' + $wNs + '/>'
            $p.Range.InsertXML($xml)
            break
        }
    }
}

# ---------------------------------------------------------------------------
# 3) Append five new paragraphs at the very end of the document body (just
#    before the sectPr), covering 5.2.1 / 5.2.2 / 5.23 / a blank line / the
#    closing non-convex remark (with a fresh "_GoBack" bookmark inside it).
# ---------------------------------------------------------------------------
$last = $d.Paragraphs.Last
$last.Range.InsertParagraphAfter()

$newLast = $d.Paragraphs.Last
$newXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><w:r><w:t>5.2</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>.</w:t></w:r><w:r><w:t>1 Backprogapation  --unrolling parameters</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">5.2.2  Gradient checking    </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>微积分定义检查得到参数是否相似</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><m:oMath><m:r><m:rPr><m:sty m:val="p"/></m:rPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>ε=</m:t></m:r><m:sSup><m:sSupPr><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr></m:ctrlPr></m:sSupPr><m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>10</m:t></m:r></m:e><m:sup><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>-4</m:t></m:r></m:sup></m:sSup></m:oMath><w:r><w:t xml:space="preserve">, </w:t></w:r><m:oMath><m:r><m:rPr><m:sty m:val="p"/></m:rPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>d</m:t></m:r><w:bookmarkStart w:id="10" w:name="OLE_LINK11"/><w:bookmarkStart w:id="11" w:name="OLE_LINK12"/><m:r><m:rPr><m:sty m:val="p"/></m:rPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>θ</m:t></m:r><w:bookmarkEnd w:id="10"/><w:bookmarkEnd w:id="11"/><m:r><m:rPr><m:sty m:val="p"/></m:rPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>=</m:t></m:r><m:f><m:fPr><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr></m:ctrlPr></m:fPr><m:num><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>J</m:t></m:r><m:d><m:dPr><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:i/></w:rPr></m:ctrlPr></m:dPr><m:e><m:r><m:rPr><m:sty m:val="p"/></m:rPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>θ</m:t></m:r><m:r><m:rPr><m:sty m:val="p"/></m:rPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>+e</m:t></m:r></m:e></m:d><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>-J(</m:t></m:r><m:r><m:rPr><m:sty m:val="p"/></m:rPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>θ</m:t></m:r><m:r><m:rPr><m:sty m:val="p"/></m:rPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>-e</m:t></m:r><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>)</m:t></m:r></m:num><m:den><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>2e</m:t></m:r></m:den></m:f></m:oMath><w:r><w:t xml:space="preserve">   </w:t></w:r></w:p><w:p><w:r><w:t>5.23   Random initialization  (</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>打破对称，消除冗余计算</w:t></w:r><w:r><w:t>)</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">            </w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">          </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>NN</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">　</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">non-convex </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>（非凸优化</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>，只能找到局部最小值</w:t></w:r><w:bookmarkStart w:id="12" w:name="_GoBack"/><w:bookmarkEnd w:id="12"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>）</w:t></w:r></w:p>
'@

$newLast.Range.InsertXML($newXml)

Write-Output "ok"
